# Fruta / hortaliza, semanal
# Weekly update: a new daily price record for "Mango" at Vega Monumental
# Concepción is inserted as row 69 (pushing the previously existing rows
# 69-144 down to 70-145), carrying the same product/category metadata as
# the old row 69 but with a new date and a different country of origin.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 69 (shifts rows 69:144 -> 70:145).
$ws.Rows.Item(69).Insert()

# Populate the newly inserted row 69 with the new weekly record.
$ws.Cells.Item(69, 1).Value  = 11
$ws.Cells.Item(69, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(69, 3).Value  = "Bíobío"
$ws.Cells.Item(69, 4).Value  = 44897
$ws.Cells.Item(69, 5).Value  = 8
$ws.Cells.Item(69, 6).Value  = "Fruta"
$ws.Cells.Item(69, 7).Value  = 100108
$ws.Cells.Item(69, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(69, 9).Value  = 100108002
$ws.Cells.Item(69, 10).Value = "Mango"
$ws.Cells.Item(69, 11).Value = "Sin especificar"
$ws.Cells.Item(69, 12).Value = "Primera"
$ws.Cells.Item(69, 13).Value = 200
$ws.Cells.Item(69, 14).Value = 8000
$ws.Cells.Item(69, 15).Value = 8500
$ws.Cells.Item(69, 16).Value = 8250
$ws.Cells.Item(69, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(69, 18).Value = "Perú"
$ws.Cells.Item(69, 19).Value = 2062
$ws.Cells.Item(69, 20).Value = 4
